$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells we are about to rewrite to stay text
# (the source data stores prices as literal strings, e.g. "1.000", "26.442.64")
# so Excel does not silently convert number-looking text into real numbers.
# (Row 44's price is untouched by this update, so it is deliberately excluded.)
$ws.Range('D2:D43').NumberFormat = '@'
$ws.Range('D45:D51').NumberFormat = '@'

$ws.Range('D2').Value = '26.442.64'
$ws.Range('E2').Value = '  -0.62%  '

$ws.Range('D3').Value = '1.842.50'
$ws.Range('E3').Value = '  -0.90%  '

$ws.Range('D4').Value = '0.9991'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('D5').Value = '261.93'
$ws.Range('E5').Value = '  -4.06%  '

$ws.Range('D6').Value = '0.9997'
$ws.Range('E6').Value = '  +0.00%  '

$ws.Range('D7').Value = '0.5190'
$ws.Range('E7').Value = '  -1.59%  '

$ws.Range('D8').Value = '0.3269'
$ws.Range('E8').Value = '  -2.98%  '

$ws.Range('D9').Value = '0.06793'
$ws.Range('E9').Value = '  -0.03%  '

$ws.Range('D10').Value = '18.70'
$ws.Range('E10').Value = '  -5.87%  '

$ws.Range('D11').Value = '0.7791'
$ws.Range('E11').Value = '  -1.85%  '

$ws.Range('D12').Value = '0.07758'
$ws.Range('E12').Value = '  +0.31%  '

$ws.Range('D13').Value = '1.841.80'
$ws.Range('E13').Value = '  -2.58%  '

$ws.Range('D14').Value = '87.78'
$ws.Range('E14').Value = '  -2.11%  '

$ws.Range('D15').Value = '5.003'
$ws.Range('E15').Value = '  -2.48%  '

$ws.Range('D16').Value = '0.9991'
$ws.Range('E16').Value = '  -0.01%  '

$ws.Range('D17').Value = '13.92'
$ws.Range('E17').Value = '  -3.35%  '

$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '0.000007983'
$ws.Range('E18').Value = '  -0.06%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').Value = '  +0.00%  '

$ws.Range('D20').Value = '26.460.33'
$ws.Range('E20').Value = '  -0.67%  '

$ws.Range('D21').Value = '2.069.59'
$ws.Range('E21').Value = '  -2.59%  '

$ws.Range('D22').Value = '4.617'
$ws.Range('E22').Value = '  -2.36%  '

$ws.Range('D23').Value = '9.538'
$ws.Range('E23').Value = '  -4.46%  '

$ws.Range('D24').Value = '5.976'
$ws.Range('E24').Value = '  -2.17%  '

$ws.Range('D25').Value = '144.84'
$ws.Range('E25').Value = '  -0.59%  '

$ws.Range('D26').Value = '2.187'
$ws.Range('E26').Value = '  -7.16%  '

$ws.Range('D27').Value = '1.648'
$ws.Range('E27').Value = '  -0.41%  '

$ws.Range('D28').Value = '16.94'
$ws.Range('E28').Value = '  -1.19%  '

$ws.Range('D29').Value = '111.78'
$ws.Range('E29').Value = '  -0.13%  '

$ws.Range('D30').Value = '4.165'
$ws.Range('E30').Value = '  -3.18%  '

$ws.Range('D31').Value = '4.117'
$ws.Range('E31').Value = '  -4.07%  '

$ws.Range('D32').Value = '0.08695'
$ws.Range('E32').Value = '  -2.15%  '

$ws.Range('D33').Value = '0.04828'
$ws.Range('E33').Value = '  -1.67%  '

$ws.Range('D34').Value = '0.7247'
$ws.Range('E34').Value = '  -0.45%  '

$ws.Range('D35').Value = '1.129'
$ws.Range('E35').Value = '  -2.45%  '

$ws.Range('D36').Value = '2.845'
$ws.Range('E36').Value = '  -1.16%  '

$ws.Range('D37').Value = '3.091'
$ws.Range('E37').Value = '  -4.05%  '

$ws.Range('D38').Value = '0.01779'
$ws.Range('E38').Value = '  -3.26%  '

$ws.Range('D39').Value = '2.225'
$ws.Range('E39').Value = '  -4.40%  '

$ws.Range('D40').Value = '0.4847'
$ws.Range('E40').Value = '  -4.74%  '

$ws.Range('D41').Value = '0.9069'
$ws.Range('E41').Value = '  -3.29%  '

$ws.Range('D42').Value = '111.33'
$ws.Range('E42').Value = '  -4.00%  '

$ws.Range('D43').Value = '6.064'
$ws.Range('E43').Value = '  -1.18%  '

$ws.Range('E44').Value = '  +0.05%  '

$ws.Range('D45').Value = '7.729'
$ws.Range('E45').Value = '  -3.35%  '

$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '0.4173'
$ws.Range('E46').Value = '  -5.41%  '

$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05926'
$ws.Range('E47').Value = '  -0.08%  '

$ws.Range('D48').Value = '9.023'
$ws.Range('E48').Value = '  -2.90%  '

$ws.Range('D49').Value = '35.03'
$ws.Range('E49').Value = '  -2.74%  '

$ws.Range('D50').Value = '0.1229'
$ws.Range('E50').Value = '  -7.16%  '

$ws.Range('D51').Value = '0.8856'
$ws.Range('E51').Value = '  +0.91%  '
